$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update realeffort (re-simulated values), and re-sort a few adjacent rows
# within race groups whose order flips under the new realeffort values
# (index/prolificid/name/gender move together with the row's rank).
$ws.Range("H2").Value = 11.35146450363736
$ws.Range("H3").Value = 10.46509312749224
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "5f2c1a97a6809c060fec8820"
$ws.Range("F4").Value = "Maggie"
$ws.Range("H4").Value = 8.345780257992518
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = "60a71d27a66fac796ad4de6f"
$ws.Range("F5").Value = "Jennifer"
$ws.Range("H5").Value = 8.335815743434177
$ws.Range("H6").Value = 7.149319585641714
$ws.Range("H7").Value = 6.428054772178267
$ws.Range("H8").Value = 6.387792560901622
$ws.Range("H9").Value = 5.358859522459975
$ws.Range("H10").Value = 5.046880051605173
$ws.Range("H11").Value = 4.123973324417245
$ws.Range("H12").Value = 2.400210451344222
$ws.Range("H13").Value = 0.1240641252646651
$ws.Range("H14").Value = 8.43712372573634
$ws.Range("H15").Value = 8.085171174491482
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "60ba8ba51a5e0a105396888a"
$ws.Range("F16").Value = "Alfredo"
$ws.Range("G16").Value = "male"
$ws.Range("H16").Value = 7.483015296297952
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("F17").Value = "Valeria"
$ws.Range("G17").Value = "female"
$ws.Range("H17").Value = 7.18543091573438
$ws.Range("H18").Value = 6.391489871176138
$ws.Range("H19").Value = 6.387011644639443
$ws.Range("H20").Value = 5.012348235563821
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = "5e706891c396cc64388ef760"
$ws.Range("F21").Value = "Maria"
$ws.Range("H21").Value = 3.142757349846526
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = "5e0adc8f4cac6834756db412"
$ws.Range("F22").Value = "Mary"
$ws.Range("H22").Value = 3.037402404511541
$ws.Range("H23").Value = 2.344227293246886
$ws.Range("H24").Value = 1.337899327810482
$ws.Range("H25").Value = 0.002543171126171584
